$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Дата добавления" (date added) note column so every row shows
# the same value "BSA 00, 13.01.25" instead of the older mixed values.
$ws.Range("E3").Value = "BSA 00, 13.01.25"
$ws.Range("E4").Value = "BSA 00, 13.01.25"

# Update the view/selection state to match the authored file
$ws.Range("E4").Select()
